$d = $word.ActiveDocument

# 1. Change highlight color from green to yellow on every paragraph that
#    actually has highlighted text (skip the trailing empty paragraphs).
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.Trim().Length -gt 0) {
        $p.Range.Font.HighlightColorIndex = 7   # wdYellow
    }
}

# 2. Remove the stale "_GoBack" bookmark (left over from the previous edit
#    session); the new edit session leaves its own trail of "_Hlk..."
#    bookmarks instead.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3. Re-create the "_Hlk..." navigation bookmarks that Word drops in as a
#    side effect of selecting/editing each of these paragraphs in turn -
#    one bookmark per paragraph, each opened at the start of its paragraph
#    and naturally closed once the next paragraph's bookmark is opened.
$names = @(
    "_Hlk208875280",
    "_Hlk208876290",
    "_Hlk208876416",
    "_Hlk208876640",
    "_Hlk208877194",
    "_Hlk208877490",
    "_Hlk208877610"
)
for ($i = 0; $i -lt $names.Length; $i++) {
    $p = $d.Paragraphs.Item(8 + $i)
    $r = $p.Range.Duplicate()
    $d.Bookmarks.Add($names[$i], $r)
}
